$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 3 becomes the "Siapa yang dapat menggunakan Aplikasi DARA?" FAQ,
#     replacing the old "Apa saja fitur utama..." entry, and loses its
#     special (taller / indented) formatting in favour of the plain look
#     used by the rest of the table.
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Siapa yang dapat menggunakan Aplikasi DARA?"
$ws.Range("C3").Value = "Manajemen rumah sakit, Instansi pemerintah yang terkait dengan layanan kesehatan, Peneliti dan akademisi yang mempelajari layanan publik"

# --- Row 4 ---
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Data apa saja yang dianalisis oleh DARA?"
$ws.Range("C4").Value = "DARA menganalisis data ulasan dari platform Google Maps, media sosial, atau sumber lainnya yang relevan. Data ini diproses secara otomatis untuk menghasilkan wawasan mengenai layanan rumah sakit."

# --- Row 5 ---
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "Bagaimana cara kerja DARA?"
$ws.Range("C5").Value = "Aplikasi DARA menggunakan metode Naive Bayes untuk analisis topik dan sentimen. Data ulasan diproses melalui model pembelajaran mesin yang telah dilatih untuk mengenali pola dalam teks, kemudian hasil analisis disajikan dalam dashboard interaktif."

# --- Row 6 ---
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "Bagaimana cara mengakses aplikasi DARA?"
$ws.Range("C6").Value = "Aplikasi DARA dapat diakses melalui browser atau perangkat mobile. Silakan hubungi tim pengembang untuk informasi lebih lanjut tentang akses dan penggunaannya."

# --- Row 7 ---
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "Apakah DARA mendukung integrasi dengan sistem lain?"

# --- Row 8 (set before row 7's answer so new shared strings are appended
#     in the same order the reference workbook used) ---
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "Bagaimana cara memberikan masukan atau melaporkan masalah terkait DARA?"
$ws.Range("C8").Value = "Anda dapat menghubungi tim dukungan pelanggan DARA melalui email: layanan.dara@gmail.com.com atau mengisi formulir kontak di situs resmi aplikasi."

$ws.Range("C7").Value = "Ya, DARA dapat diintegrasikan dengan sistem informasi rumah sakit (SIMRS) atau platform lain melalui API yang tersedia."

# --- The last two FAQ rows (old "Apa manfaat utama..." and the old
#     "...support@dara.com.com" entry) are removed entirely. ---
$ws.Rows.Item(9).Delete()
$ws.Rows.Item(9).Delete()

# --- Normalise row 3's look: drop the taller row height and the special
#     indent/alignment styling, keeping only the bold question font that
#     every other question cell uses. ---
$ws.Range("A3").Style = "Normal"
$ws.Range("C3").Style = "Normal"
$ws.Range("B3").Style = "Normal"
$ws.Range("B3").Font.Bold = $true
$ws.Rows.Item(3).AutoFit()

# --- Scroll back to the left edge and select A2:A8. ---
$ws.Range("A2:A8").Select() | Out-Null
